$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '98.809.73'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.46%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.357.74'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.36%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '257.21'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.11%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '636.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.58'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +24.87%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.413'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +7.28%  '

$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.09'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +33.87%  '

$ws.Range("B10").Value = 'USDC'
$ws.Range("C10").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.999'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.359.77'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.57%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.206'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.41%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.84'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +20.61%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '98.536.01'
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000255'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.64%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.984.29'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.33%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.49'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.11%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.352.68'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.02%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.49'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +11.05%  '

$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.85'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +16.13%  '

$ws.Range("B21").Value = 'SuiNetwork'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '3.54'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.83%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '494.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.18%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.90'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.50%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000208'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.51%  '

$ws.Range("B25").Value = 'Stellar'
$ws.Range("C25").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.395'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +44.54%  '

$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.98'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.36%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '91.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.22%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.40'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.91%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.521.80'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.42%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.153'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +22.76%  '

$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.25'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +22.96%  '

$ws.Range("B32").Value = 'Dai'
$ws.Range("C32").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.01%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.186'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.02%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '29.16'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.40%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.508'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +12.44%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.62'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.05%  '

$ws.Range("E38").Value = '  +0.12%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.02'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.16%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '515.28'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.47%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '24.73'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.32%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.86'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.88%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.29'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.66%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.823'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.27%  '

$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.24'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.40%  '

$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.09%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.03'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '162.47'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.50%  '

$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0368'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +15.69%  '

$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.58'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +18.73%  '

$ws.Range("B51").Value = 'ImmutableX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.48'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.92%  '
